# Adds a new '2022-Q1' worksheet (fund-holdings detail) before '总计',
# and prepends a summary row for it on the '总计' sheet.

function Set-TextCell($ws, $row, $col, $val) {
    # Forces the cell to be stored as text (keeps leading zeros / avoids
    # numeric auto-conversion), then strips the number-format override so
    # the cell is left with no explicit style, matching the sibling sheets.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

function Set-NumCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $val
    $c.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- 1. Create the '2022-Q1' sheet by duplicating '2021-Q4' (same column
#        layout/header/style) and placing it right before '总计'. ---
$totalSheet = $wb.Worksheets.Item("总计")
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($totalSheet)
$ns = $wb.ActiveSheet
$ns.Name = "2022-Q1"

# Header row (row 1) already matches: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名

$data2022 = @(
    @(0, "519692", "交银成长混合A", "23.39", "82.72", "8.76", "2.0490", 2),
    @(1, "960016", "交银成长混合H", "23.39", "82.72", "8.76", "2.0490", 2),
    @(2, "519772", "交银新生活力灵活配置混合", "54.07", "86.33", "3.17", "1.7140", 7),
    @(3, "519694", "交银蓝筹混合", "16.14", "82.45", "8.76", "1.4139", 2),
    @(4, "519773", "交银施罗德数据产业灵活配置混合", "19.72", "86.28", "6.14", "1.2108", 3),
    @(5, "010094", "交银施罗德产业机遇混合", "18.53", "85.59", "6.29", "1.1655", 3),
    @(6, "519732", "交银定期支付双息平衡混合", "40.83", "67.67", "2.06", "0.8411", 6),
    @(7, "010654", "天弘医药创新混合A", "11.86", "87.22", "6.68", "0.7922", 3),
    @(8, "501093", "华夏翔阳两年定期开放混合", "7.35", "78.10", "4.42", "0.3249", 5),
    @(9, "010655", "天弘医药创新混合C", "3.52", "87.22", "6.68", "0.2351", 3),
    @(10, "168207", "中融创业板两年定期开放混合", "3.85", "88.46", "5.26", "0.2025", 1),
    @(11, "001387", "中融新经济灵活配置混合A", "3.35", "93.13", "5.12", "0.1715", 2),
    @(12, "004905", "华泰柏瑞生物医药灵活配置混合A", "5.73", "87.55", "2.76", "0.1581", 10),
    @(13, "001701", "中融产业升级灵活配置混合", "2.84", "92.91", "5.23", "0.1485", 2),
    @(14, "005805", "华泰柏瑞医疗健康混合A", "4.88", "88.16", "2.81", "0.1371", 9),
    @(15, "004784", "招商稳健优选股票", "4.90", "85.49", "2.71", "0.1328", 8),
    @(16, "010697", "中融行业先锋6个月持有期混合A", "1.79", "89.23", "5.05", "0.0904", 2),
    @(17, "001388", "中融新经济灵活配置混合C", "0.72", "93.13", "5.12", "0.0369", 2),
    @(18, "010613", "中融产业趋势一年定期开放混合A", "0.69", "87.17", "4.38", "0.0302", 9),
    @(19, "007775", "汇安量化先锋混合A", "0.38", "94.51", "3.04", "0.0116", 9),
    @(20, "006240", "中融医疗健康精选混合A", "0.14", "94.05", "5.00", "0.0070", 4),
    @(21, "010031", "华泰柏瑞生物医药灵活配置混合C", "0.25", "87.55", "2.76", "0.0069", 10),
    @(22, "010614", "中融产业趋势一年定期开放混合C", "0.11", "87.17", "4.38", "0.0048", 9),
    @(23, "010698", "中融行业先锋6个月持有期混合C", "0.09", "89.23", "5.05", "0.0045", 2),
    @(24, "006241", "中融医疗健康精选混合C", "0.08", "94.05", "5.00", "0.0040", 4),
    @(25, "011453", "华泰柏瑞医疗健康混合C", "0.14", "88.16", "2.81", "0.0039", 9),
    @(26, "007776", "汇安量化先锋混合C", "0.11", "94.51", "3.04", "0.0033", 9)
)

for ($i = 0; $i -lt $data2022.Count; $i++) {
    $row = $data2022[$i]
    $r = $i + 2
    Set-NumCell  $ns $r 1 $row[0]   # A: row index
    Set-TextCell $ns $r 2 $row[1]   # B: 基金代码
    Set-TextCell $ns $r 3 $row[2]   # C: 基金名称
    Set-TextCell $ns $r 4 $row[3]   # D: 基金规模
    Set-TextCell $ns $r 5 $row[4]   # E: 股票总仓位
    Set-TextCell $ns $r 6 $row[5]   # F: 仓位占比
    Set-TextCell $ns $r 7 $row[6]   # G: 持有市值(亿元)
    Set-NumCell  $ns $r 8 $row[7]   # H: 仓位排名
}

# Give the A column (row index) the same bold/bordered style used by the
# header row & other quarter sheets, matching rows up to the new max row.
for ($r = 2; $r -le ($data2022.Count + 1); $r++) {
    $ns.Cells.Item(2, 1).Copy($ns.Cells.Item($r, 1))
    $ns.Cells.Item($r, 1).Value = $data2022[$r - 2][0]
}

# --- 2. Insert a new top data row on '总计' summarising 2022-Q1. ---
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()

# Re-use the row-index column's style (bold/bordered, s=2) for the new row.
$tot.Cells.Item(3,1).Copy($tot.Cells.Item(2,1))
$tot.Cells.Item(2,1).Value = 0

Set-TextCell $tot 2 2 "2022-Q1"   # 日期
Set-NumCell  $tot 2 3 27   # 持有数量(只)
Set-NumCell  $tot 2 4 12.95   # 持有市值(亿元)

# Renumber the row-index column for the rows pushed down.
for ($r = 3; $r -le 5; $r++) {
    $tot.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "2022-Q1 sheet added and 总计 updated"
